$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 2988.6667
$ws.Range("I4").Value = 2988.6667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2988.6667
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = -2874.6667
$ws.Range("M4").ClearContents()

# Row 10
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = 5
$ws.Range("N10").Value = -591
$ws.Range("L10").ClearContents()

# Row 28
$ws.Range("H28").Value = 864.8570999999999
$ws.Range("I28").Value = 809.6
$ws.Range("K28").Value = 809.6
$ws.Range("M28").Value = -324.6

# Row 40
$ws.Range("H40").Value = 8666.666999999999
$ws.Range("J40").Value = 8666.666999999999
$ws.Range("L40").Value = 8666.666999999999
$ws.Range("N40").Value = -9016.666999999999

# Row 70
$ws.Range("H70").Value = 21350.2
$ws.Range("J70").Value = 35083
$ws.Range("L70").Value = 105249
$ws.Range("N70").Value = -105789

# Row 73
$ws.Range("H73").Value = 21350.2
$ws.Range("J73").Value = 35083
$ws.Range("L73").Value = 105249
$ws.Range("N73").Value = -107121

# Row 132
$ws.Range("H132").Value = 3807.4666
$ws.Range("J132").Value = 3054.3333
$ws.Range("L132").Value = 9162.999899999999
$ws.Range("N132").Value = -14222.9999

# Row 134
$ws.Range("H134").Value = 95000
$ws.Range("J134").Value = 95000
$ws.Range("L134").Value = 95000
$ws.Range("N134").Value = -105140

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2899.5
$ws.Range("I61").Value = 2899.5
$ws.Range("K61").Value = 2899.5
$ws.Range("M61").Value = -2687.5

# Row 88
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("N88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").ClearContents()

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("N91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").ClearContents()

# Row 136
$ws.Range("H136").Value = 2899.5
$ws.Range("I136").Value = 2899.5
$ws.Range("K136").Value = 8698.5
$ws.Range("M136").Value = -6148.5

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 4500.6665
$ws.Range("I22").Value = 5500
$ws.Range("J22").Value = 2502
$ws.Range("K22").Value = 5500
$ws.Range("L22").Value = 2502
$ws.Range("M22").Value = -5327
$ws.Range("N22").Value = -2848

# Row 86
$ws.Range("H86").Value = 1324.625
$ws.Range("I86").Value = 1282.8334
$ws.Range("J86").Value = 1450
$ws.Range("K86").Value = 1282.8334
$ws.Range("L86").Value = 1450
$ws.Range("M86").Value = -159.8334
$ws.Range("N86").Value = -3696

# Row 89
$ws.Range("H89").Value = 1324.625
$ws.Range("I89").Value = 1282.8334
$ws.Range("J89").Value = 1450
$ws.Range("K89").Value = 6414.166999999999
$ws.Range("L89").Value = 7250
$ws.Range("M89").Value = -798.1669999999995
$ws.Range("N89").Value = -18482

# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("N101").Value = 0
$ws.Range("L101").ClearContents()

# Row 110
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").Value = 0
$ws.Range("L132").ClearContents()

# Row 134
$ws.Range("H134").Value = 4858.909
$ws.Range("I134").Value = 3778.4285
$ws.Range("J134").Value = 6749.75
$ws.Range("K134").Value = 11335.2855
$ws.Range("L134").Value = 20249.25
$ws.Range("M134").Value = -8800.2855
$ws.Range("N134").Value = -25319.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1016.6667
$ws.Range("I16").Value = 1016.6667
$ws.Range("K16").Value = 1016.6667
$ws.Range("M16").Value = -729.6667

# Row 58
$ws.Range("H58").Value = 1099.8
$ws.Range("J58").Value = 1249.5
$ws.Range("L58").Value = 1249.5
$ws.Range("N58").Value = -1655.5

# Row 69
$ws.Range("H69").Value = 11266.333
$ws.Range("I69").Value = 6899.5
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 6899.5
$ws.Range("L69").Value = 20000
$ws.Range("M69").Value = -6150.5
$ws.Range("N69").Value = -21498

# Row 72
$ws.Range("H72").Value = 11266.333
$ws.Range("I72").Value = 6899.5
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 20698.5
$ws.Range("L72").Value = 60000
$ws.Range("M72").Value = -16954.5
$ws.Range("N72").Value = -67488

# Row 86
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("N86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()

# Row 89
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("N89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()

# Row 99
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("L100").ClearContents()

# Row 113
$ws.Range("H113").Value = 1016.6667
$ws.Range("I113").Value = 1016.6667
$ws.Range("K113").Value = 1016.6667
$ws.Range("M113").Value = 1153.3333

# Row 126
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

# Row 132
$ws.Range("H132").Value = 1847
$ws.Range("I132").Value = 1931.7
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5795.1
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3265.1
$ws.Range("N132").Value = -8060

# Row 136
$ws.Range("H136").Value = 1099.8
$ws.Range("J136").Value = 1249.5
$ws.Range("L136").Value = 3748.5
$ws.Range("N136").Value = -8848.5

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 2504000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2504000
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 2504000
$ws.Range("N7").Value = -2504224
$ws.Range("L7").ClearContents()

# Row 8
$ws.Range("H8").Value = 2504000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2504000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2504000
$ws.Range("M8").Value = -5361.5
$ws.Range("N8").Value = -2504278

# Row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

# Row 55
$ws.Range("H55").Value = 60000
$ws.Range("J55").Value = 60000
$ws.Range("L55").Value = 60000
$ws.Range("N55").Value = -60654

# Row 133
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120

# Row 135
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 132
$ws.Range("H132").Value = 2202.6
$ws.Range("I132").Value = 2202.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6607.799999999999
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -4077.799999999999
$ws.Range("M132").ClearContents()

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("N138").Value = 0
$ws.Range("L138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 10006
$ws.Range("I15").Value = 10006
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 10006
$ws.Range("N15").Value = 0
$ws.Range("M15").Value = -9718
$ws.Range("L15").ClearContents()

# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").Value = 0
$ws.Range("L75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").Value = 0
$ws.Range("L78").ClearContents()

# Row 122
$ws.Range("H122").Value = 2123.9333
$ws.Range("I122").Value = 2127.2307
$ws.Range("K122").Value = 6381.6921
$ws.Range("M122").Value = -3931.6921

# Row 132
$ws.Range("H132").Value = 5549.9287
$ws.Range("I132").Value = 3979.9
$ws.Range("J132").Value = 9475
$ws.Range("K132").Value = 11939.7
$ws.Range("L132").Value = 28425
$ws.Range("M132").Value = -9409.700000000001
